# "issues fixes and scenarios implementation"
# Populate Sheet1 with the sales-report header row (A1:O1), wrap the header
# text, give the "Sales to taxable person (Value)" header vertical centering
# too, grow row 1 to fit the wrapped text, and finish with O1 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$headers = @(
    "Account Code",
    "Account Description",
    "Account Reference",
    "Signature Date",
    "Branch Name",
    "Tax Period",
    "Invoice Date",
    "Invoice Number",
    "Description",
    "Quantity",
    "Non Taxable Sales",
    "Value Of Exports",
    "Sales to taxable person (Value)",
    "Sales to Consumer (Value)",
    "Item subject to taxes"
)

for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.WrapText = $true
}

# "Sales to taxable person (Value)" (column M) also gets vertical centering.
$ws.Cells.Item(1, 13).VerticalAlignment = -4108

# Header row is tall enough to show the wrapped text.
$ws.Rows.Item(1).RowHeight = 60

$ws.Range("O1").Select() | Out-Null
